# Update cryptocurrency price list (GitHub Actions scheduled refresh)
# Applies the latest scraped prices / 1h volume percentages onto Sheet1,
# and replaces the former last row (EnergySwap) with BabyDogeCoin while
# pushing Algorand down one row, matching the refreshed coinranking feed.
#
# Price/volume columns hold text (not numbers) in the source sheet, so
# numeric-looking values are written with a leading apostrophe - same as
# typing them straight into Excel - to keep them stored as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.694.89'
$ws.Range("E2").Value = '  +1.11%  '
$ws.Range("D3").Value = '1.644.55'
$ws.Range("E3").Value = '  +0.34%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '''213.09'
$ws.Range("E5").Value = '  +0.82%  '
$ws.Range("E6").Value = '  -0.60%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '''23.35'
$ws.Range("E8").Value = '  +1.29%  '
$ws.Range("D9").Value = '''0.259'
$ws.Range("E9").Value = '  +1.08%  '
$ws.Range("E11").Value = '  +0.53%  '
$ws.Range("D12").Value = '1.877.94'
$ws.Range("D13").Value = '1.637.75'
$ws.Range("E13").Value = '  -0.27%  '
$ws.Range("E14").Value = '  +0.62%  '
$ws.Range("D15").Value = '''0.563'
$ws.Range("E15").Value = '  +0.81%  '
$ws.Range("D16").Value = '''64.73'
$ws.Range("E16").Value = '  +0.68%  '
$ws.Range("D17").Value = '27.688.48'
$ws.Range("E17").Value = '  +1.19%  '
$ws.Range("D18").Value = '''231.55'
$ws.Range("E19").Value = '  +1.98%  '
$ws.Range("E20").Value = '  +0.54%  '
$ws.Range("E21").Value = '  +0.01%  '
$ws.Range("D22").Value = '''4.29'
$ws.Range("E22").Value = '  -0.59%  '
$ws.Range("D23").Value = '''10.04'
$ws.Range("E23").Value = '  +5.70%  '
$ws.Range("E24").Value = '  -2.98%  '
$ws.Range("D25").Value = '''150.05'
$ws.Range("E25").Value = '  +1.77%  '
$ws.Range("D26").Value = '''6.95'
$ws.Range("E26").Value = '  -0.17%  '
$ws.Range("E27").Value = '  -1.21%  '
$ws.Range("E28").Value = '  +0.99%  '
$ws.Range("E29").Value = '  +0.07%  '
$ws.Range("E30").Value = '  +0.95%  '
$ws.Range("E31").Value = '  +0.71%  '
$ws.Range("D32").Value = '''3.31'
$ws.Range("E32").Value = '  +1.24%  '
$ws.Range("D33").Value = '1.455.19'
$ws.Range("E33").Value = '  +3.35%  '
$ws.Range("E34").Value = '  +0.75%  '
$ws.Range("E35").Value = '  +0.66%  '
$ws.Range("E36").Value = '  -1.06%  '
$ws.Range("E37").Value = '  +1.10%  '
$ws.Range("D38").Value = '''0.884'
$ws.Range("E38").Value = '  +0.49%  '
$ws.Range("E39").Value = '  +0.51%  '
$ws.Range("D40").Value = '''0.885'
$ws.Range("E40").Value = '  +11.76%  '
$ws.Range("D41").Value = '''71.08'
$ws.Range("E41").Value = '  +10.23%  '
$ws.Range("E42").Value = '  +0.38%  '
$ws.Range("E43").Value = '  +0.05%  '
$ws.Range("E44").Value = '  +2.22%  '
$ws.Range("E45").Value = '  +0.80%  '
$ws.Range("D47").Value = '1.787.45'
$ws.Range("E47").Value = '  +0.41%  '
$ws.Range("E48").Value = '  +5.17%  '
$ws.Range("D49").Value = '''85.76'
$ws.Range("E49").Value = '  -1.73%  '
$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").Value = '0.0₆0107'
$ws.Range("E50").Value = '  +1.73%  '
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").Value = '''0.0990'
$ws.Range("E51").Value = '  +0.16%  '
